$d = $word.ActiveDocument

# --- Change 1: first paragraph -------------------------------------------
# "This is a Microsoft word document." -> same text + two trailing spaces,
# followed by a new red run: "(This is a change – Version for branch alternate)"

$p1 = $d.Paragraphs(1).Range
$d.Content.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

$p1 = $d.Paragraphs(1).Range
$insertStart = $p1.End
$newText = "(This is a change " + [char]0x2013 + " Version for branch alternate)"
$p1.InsertAfter($newText)
$newRange = $d.Range($insertStart - 1, $p1.End - 1)
$newRange.Font.Color = 192

# --- Change 2: last "Bop-bop-bop-bop" run ---------------------------------
# Split it into its own run (bracketed by proofErr spellStart/spellEnd),
# matching the other "Bop-bop-bop-bop" occurrences in the document, while
# keeping its run formatting identical.

$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$target = $d.Range($lastPara.End - 16, $lastPara.End - 1)
$target.Bold = $target.Bold
